$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "26.974.60"

$ws.Range("D3").Value = "1.561.18"
$ws.Range("E3").Value = "  +0.38%  "

Set-TextValue $ws.Range("D4") "1.01"
$ws.Range("E4").Value = "  +0.23%  "

Set-TextValue $ws.Range("D5") "207.36"
$ws.Range("E5").Value = "  +0.44%  "

$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("E8").Value = "  +0.88%  "

$ws.Range("E9").Value = "  +0.64%  "

Set-TextValue $ws.Range("D10") "0.0601"
$ws.Range("E10").Value = "  +2.56%  "

$ws.Range("E11").Value = "  -0.39%  "

$ws.Range("D12").Value = "1.783.96"
$ws.Range("E12").Value = "  +0.41%  "

$ws.Range("D13").Value = "1.567.08"
$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("E14").Value = "  +0.45%  "

$ws.Range("E15").Value = "  +0.65%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "26.961.67"
$ws.Range("E16").Value = "  +0.15%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D17") "61.88"
$ws.Range("E17").Value = "  +0.24%  "

Set-TextValue $ws.Range("D19") "215.76"
$ws.Range("E19").Value = "  -1.20%  "

$ws.Range("E20").Value = "  +1.05%  "

$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("E22").Value = "  +1.58%  "

$ws.Range("E23").Value = "  +0.07%  "

Set-TextValue $ws.Range("D24") "1.92"
$ws.Range("E24").Value = "  -0.85%  "

Set-TextValue $ws.Range("D25") "153.45"
$ws.Range("E25").Value = "  -0.32%  "

Set-TextValue $ws.Range("D27") "15.08"
$ws.Range("E27").Value = "  +1.15%  "

$ws.Range("E28").Value = "  +1.47%  "

$ws.Range("E29").Value = "  +0.07%  "

Set-TextValue $ws.Range("D31") "1.11"
$ws.Range("E31").Value = "  +1.39%  "

Set-TextValue $ws.Range("D32") "3.22"
$ws.Range("E32").Value = "  +0.29%  "

$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").Value = "1.423.59"
$ws.Range("E33").Value = "  -0.66%  "

$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D34") "3.11"
$ws.Range("E34").Value = "  +1.57%  "

$ws.Range("E35").Value = "  +2.54%  "

$ws.Range("E36").Value = "  +8.05%  "

Set-TextValue $ws.Range("D37") "2.35"
$ws.Range("E37").Value = "  +2.46%  "

$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("E39").Value = "  +1.98%  "

$ws.Range("E41").Value = "  -0.27%  "

$ws.Range("E42").Value = "  +0.09%  "

$ws.Range("E43").Value = "  +2.71%  "

$ws.Range("E44").Value = "  +1.98%  "

Set-TextValue $ws.Range("D45") "64.53"
$ws.Range("E45").Value = "  +1.09%  "

$ws.Range("E46").Value = "  -0.44%  "

$ws.Range("D47").Value = "1.697.99"
$ws.Range("E47").Value = "  +0.42%  "

Set-TextValue $ws.Range("D48") "87.19"
$ws.Range("E48").Value = "  +0.32%  "

$ws.Range("E49").Value = "  -0.53%  "

$ws.Range("E50").Value = "  +0.04%  "

Set-TextValue $ws.Range("D51") "0.0958"
$ws.Range("E51").Value = "  +0.62%  "
